# Update countries & provincias Spain
#
# Daily refresh of the COVID-19 "Pais" worksheet:
#   - updated case counters for several countries
#   - "Burkina Faso" overtook "Cuba" / "Reunion" / "Jordania" in total
#     cases, so those four rows re-sort (Burkina Faso jumps to the top
#     of that block, the other three each drop one place)
#   - "Bermudas" overtook "San Martin (Parte Holandesa)", so those two
#     rows swap
#   - refreshed the "Datos actualizados..." timestamp cell (A1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Country, $B, $C, $D, $E, $F, $G, $H) {
    $ws.Range("A$Row").Value = $Country
    $ws.Range("B$Row").Value = $B
    $ws.Range("C$Row").Value = $C
    $ws.Range("D$Row").Value = $D
    $ws.Range("E$Row").Value = $E
    $ws.Range("F$Row").Value = $F
    $ws.Range("G$Row").Value = $G
    $ws.Range("H$Row").Value = $H
}

# Estados Unidos
Set-Row 4 "Estados Unidos" 363408 26735 19313 333332 8879 1147 10763

# Austria
Set-Row 17 "Austria" 12297 246 3463 8614 250 16 220

# Finlandia
Set-Row 44 "Finlandia" 2176 249 300 1849 81 0 27

# Burkina Faso climbs past Cuba, Reunion and Jordania
Set-Row 92 "Burkina Faso" 364 19 108 238 0 1 18
Set-Row 93 "Cuba" 350 30 18 323 12 1 9
Set-Row 94 "Reunion" 349 5 40 309 4 0 0
Set-Row 95 "Jordania" 349 4 126 217 5 1 6

# Bermudas overtakes San Martin (Parte Holandesa)
Set-Row 149 "Bermudas" 39 2 17 20 0 2 2
Set-Row 150 "San Martin (Parte Holandesa)" 37 12 1 30 0 2 6

# Refresh the "last updated" timestamp cell
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 23:52"
